$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$row2 = New-Object 'object[,]' 1,16
$row2[0,0] = 1
$row2[0,1] = 0.3333333333333333
$row2[0,2] = 0.009651333333333333
$row2[0,3] = 0.028954
$row2[0,4] = 0.8567539576860482
$row2[0,5] = 0.8567539576860482
$row2[0,6] = 3
$row2[0,7] = 1
$row2[0,8] = 5.828378333333333
$row2[0,9] = 17.485135
$row2[0,10] = 0.1413867973615592
$row2[0,11] = 0.1413867973615592
$row2[0,12] = 0.05625162208777777
$row2[0,13] = 0.50626459879
$row2[0,14] = 0.1211336982040711
$row2[0,15] = 0.1211336982040711
$ws.Range("E2:T2").Value = $row2

# Row 3
$row3 = New-Object 'object[,]' 1,16
$row3[0,0] = 1
$row3[0,1] = 0.3333333333333333
$row3[0,2] = 0.009651333333333333
$row3[0,3] = 0.028954
$row3[0,4] = 0.8567539576860482
$row3[0,5] = 0.8567539576860482
$row3[0,6] = 3
$row3[0,7] = 1
$row3[0,8] = 24.92162533333333
$row3[0,9] = 74.764876
$row3[0,10] = 0.604557320991465
$row3[0,11] = 0.604557320991465
$row3[0,12] = 0.2405269133004445
$row3[0,13] = 2.164742219704
$row3[0,14] = 0.5179568774075123
$row3[0,15] = 0.5179568774075123
$ws.Range("E3:T3").Value = $row3

# Row 4
$row4 = New-Object 'object[,]' 1,16
$row4[0,0] = 1
$row4[0,1] = 0.3333333333333333
$row4[0,2] = 0.009651333333333333
$row4[0,3] = 0.028954
$row4[0,4] = 0.8567539576860482
$row4[0,5] = 0.8567539576860482
$row4[0,6] = 3
$row4[0,7] = 1
$row4[0,8] = 10.47292833333333
$row4[0,9] = 31.418785
$row4[0,10] = 0.2540558816469758
$row4[0,11] = 0.2540558816469758
$row4[0,12] = 0.1010777223211111
$row4[0,13] = 0.90969950089
$row4[0,14] = 0.2176633820744648
$row4[0,15] = 0.2176633820744648
$ws.Range("E4:T4").Value = $row4

# Row 5
$row5 = New-Object 'object[,]' 1,16
$row5[0,0] = 1
$row5[0,1] = 0.3333333333333333
$row5[0,2] = 0.001613666666666667
$row5[0,3] = 0.004841
$row5[0,4] = 0.1432460423139518
$row5[0,5] = 0.1432460423139518
$row5[0,6] = 3
$row5[0,7] = 1
$row5[0,8] = 5.828378333333333
$row5[0,9] = 17.485135
$row5[0,10] = 0.1413867973615592
$row5[0,11] = 0.1413867973615592
$row5[0,12] = 0.009405059837222222
$row5[0,13] = 0.084645538535
$row5[0,14] = 0.02025309915748803
$row5[0,15] = 0.02025309915748803
$ws.Range("E5:T5").Value = $row5

# Row 6
$row6 = New-Object 'object[,]' 1,16
$row6[0,0] = 1
$row6[0,1] = 0.3333333333333333
$row6[0,2] = 0.001613666666666667
$row6[0,3] = 0.004841
$row6[0,4] = 0.1432460423139518
$row6[0,5] = 0.1432460423139518
$row6[0,6] = 3
$row6[0,7] = 1
$row6[0,8] = 24.92162533333333
$row6[0,9] = 74.764876
$row6[0,10] = 0.604557320991465
$row6[0,11] = 0.604557320991465
$row6[0,12] = 0.04021519607955556
$row6[0,13] = 0.361936764716
$row6[0,14] = 0.08660044358395272
$row6[0,15] = 0.08660044358395272
$ws.Range("E6:T6").Value = $row6

# Row 7
$row7 = New-Object 'object[,]' 1,16
$row7[0,0] = 1
$row7[0,1] = 0.3333333333333333
$row7[0,2] = 0.001613666666666667
$row7[0,3] = 0.004841
$row7[0,4] = 0.1432460423139518
$row7[0,5] = 0.1432460423139518
$row7[0,6] = 3
$row7[0,7] = 1
$row7[0,8] = 10.47292833333333
$row7[0,9] = 31.418785
$row7[0,10] = 0.2540558816469758
$row7[0,11] = 0.2540558816469758
$row7[0,12] = 0.01689981535388889
$row7[0,13] = 0.152098338185
$row7[0,14] = 0.03639249957251103
$row7[0,15] = 0.03639249957251103
$ws.Range("E7:T7").Value = $row7

Write-Output "Updated rows 2-7, columns E:T"